# ProtocoloClienteServidor.xlsx — update the "Salir de partida" protocol row
# to reflect the new message format (includes the outgoing player's index)
# and adjust the visible selection, per the Blackjack adaptation commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 = "Salir de partida". The notification message now also carries the
# index of the player who left, and the observations text documents that
# the server rebroadcasts that index to the rest of the players.
$ws.Range("D9").Value = '7$ID_partida/indexJugador'
$ws.Range("E9").Value = 'Este mensaje se envía cuando un cliente quiere salir de una partida para que el servidor elimine sus datos de la lista de jugadores de la partida en cuestión. Además, reenvia un mensaje con el inidce del jugador saliente al resto de jugadores para que lo tengan en cuenta.'

# Move the active selection to E10, matching where the author left off editing.
$ws.Range("E10").Select()
